# Update column G ("K") values on Sheet1 with newly computed strikeout
# counts (regenerated from source data, replacing the old "Strike#" values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G
$updates = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 3
    24 = 2
    25 = 1
    26 = 2
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 2
    34 = 2
    35 = 1
    36 = 2
    37 = 1
    39 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
